# daily auto push: 2025-10-08 09:30 UTC
# Append the day's new readings as row 79 on the (single) data sheet.
#
# Column A holds date-like strings (e.g. "2025/10/08") that must stay as
# literal text, matching every other row above it. Just assigning a
# "yyyy/mm/dd"-shaped string to .Value lets Excel's normal typed-input
# parsing kick in and silently turn it into a date serial number with a
# date number-format. To reproduce plain text we briefly force the cell to
# Text format before writing the value, then clear the formatting again so
# the cell ends up back on the sheet's default (unstyled) look, exactly
# like the surrounding historical rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/08"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "水"
$ws.Cells.Item($newRow, 3).Value = 18
$ws.Cells.Item($newRow, 4).Value = 15
